$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some "Price" values are plain decimal numbers (e.g. "211.05") stored as
# text in the workbook. Force those specific cells to Text format first so
# Excel does not auto-convert the assigned string into a numeric value.
$textCells = @('D5', 'D11', 'D17', 'D19', 'D20', 'D25', 'D26', 'D27', 'D31', 'D32', 'D36', 'D43', 'D48', 'D50')
foreach ($cell in $textCells) {
    $ws.Range($cell).NumberFormat = "@"
}

# Apply the updated price / volume figures from the latest cryptos refresh.
$ws.Range('D2').Value = '27.159.59'
$ws.Range('E2').Value = '  +0.68%  '
$ws.Range('D3').Value = '1.569.36'
$ws.Range('E3').Value = '  +0.78%  '
$ws.Range('E4').Value = '  +0.74%  '
$ws.Range('D5').Value = '211.05'
$ws.Range('E5').Value = '  +2.03%  '
$ws.Range('E6').Value = '  +0.46%  '
$ws.Range('E7').Value = '  +0.67%  '
$ws.Range('E8').Value = '  -0.43%  '
$ws.Range('E9').Value = '  +0.39%  '
$ws.Range('E10').Value = '  +0.68%  '
$ws.Range('D11').Value = '0.0867'
$ws.Range('E11').Value = '  +1.10%  '
$ws.Range('D12').Value = '1.792.09'
$ws.Range('E12').Value = '  +0.75%  '
$ws.Range('D13').Value = '1.547.09'
$ws.Range('E13').Value = '  -0.21%  '
$ws.Range('E14').Value = '  +0.62%  '
$ws.Range('E15').Value = '  -0.22%  '
$ws.Range('D16').Value = '27.116.60'
$ws.Range('E16').Value = '  +0.53%  '
$ws.Range('D17').Value = '62.19'
$ws.Range('E17').Value = '  +0.38%  '
$ws.Range('D18').Value = '0.0₃0702'
$ws.Range('E18').Value = '  -0.58%  '
$ws.Range('D19').Value = '215.75'
$ws.Range('D20').Value = '7.40'
$ws.Range('E20').Value = '  +0.89%  '
$ws.Range('E21').Value = '  +0.71%  '
$ws.Range('E22').Value = '  +1.30%  '
$ws.Range('E23').Value = '  -0.21%  '
$ws.Range('E24').Value = '  +0.73%  '
$ws.Range('D25').Value = '153.56'
$ws.Range('E25').Value = '  +0.40%  '
$ws.Range('D26').Value = '6.63'
$ws.Range('E26').Value = '  -0.13%  '
$ws.Range('D27').Value = '15.06'
$ws.Range('E27').Value = '  +0.30%  '
$ws.Range('E28').Value = '  +1.73%  '
$ws.Range('E29').Value = '  +0.60%  '
$ws.Range('E30').Value = '  +2.35%  '
$ws.Range('D31').Value = '0.0473'
$ws.Range('E31').Value = '  +0.80%  '
$ws.Range('D32').Value = '3.23'
$ws.Range('E33').Value = '  +2.42%  '
$ws.Range('D34').Value = '1.454.57'
$ws.Range('E34').Value = '  +2.46%  '
$ws.Range('E35').Value = '  +0.83%  '
$ws.Range('D36').Value = '1.61'
$ws.Range('E36').Value = '  -0.01%  '
$ws.Range('E37').Value = '  +1.77%  '
$ws.Range('E38').Value = '  +1.16%  '
$ws.Range('E39').Value = '  +0.72%  '
$ws.Range('E40').Value = '  +2.49%  '
$ws.Range('E41').Value = '  +0.17%  '
$ws.Range('E42').Value = '  +0.68%  '
$ws.Range('D43').Value = '2.36'
$ws.Range('E43').Value = '  +1.73%  '
$ws.Range('E44').Value = '  -0.50%  '
$ws.Range('E45').Value = '  -0.23%  '
$ws.Range('E46').Value = '  -0.93%  '
$ws.Range('D47').Value = '1.702.79'
$ws.Range('E47').Value = '  +0.65%  '
$ws.Range('D48').Value = '86.05'
$ws.Range('E48').Value = '  -1.68%  '
$ws.Range('E49').Value = '  +2.37%  '
$ws.Range('D50').Value = '0.0520'
$ws.Range('E50').Value = '  -0.16%  '
$ws.Range('E51').Value = '  +0.15%  '

# Restore the default "Normal" cell style on the text-forced cells so their
# formatting matches the rest of the untouched column.
foreach ($cell in $textCells) {
    $ws.Range($cell).Style = "Normal"
}
